$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F (matches the diff's new <cols> entry for column 6)
$ws.Columns.Item(6).ColumnWidth = 16

# New column C: for each data row (11-56) pull the field name out of the
# TypeScript-ish text in column F and build a "name='';" JS assignment.
# C11 is entered on its own (its own formula chain) and C12:C56 are filled
# as a second range, so they land in separate shared-formula groups -
# matching how the workbook was actually authored (fill-down from C12).
$ws.Range("C11").Formula = '=LEFT(F11,IFERROR(FIND("?",F11)-1, FIND(":",F11)-1)) & "='''';"'
$ws.Range("C12:C56").Formula = '=LEFT(F12,IFERROR(FIND("?",F12)-1, FIND(":",F12)-1)) & "='''';"'

# Match the author's final selection in the saved file.
[void]$ws.Range("C55:C56").Select()
